# Weekly update: insert this week's two new Brócoli quality-grade rows
# (Primera / Segunda) right after the last existing data row (480),
# pushing the previously-last week's rows (old 481-507) down to 483-509.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at 481 (each Insert() pushes everything
# below down by one and inherits formatting from the row above, same as
# Excel's normal "Insert Copied Cells"/"Insert Sheet Rows" behaviour).
$ws.Rows.Item(481).Insert()
$ws.Rows.Item(481).Insert()

# Row 481: Brócoli, Primera
$ws.Range("A481").Value = 6
$ws.Range("B481").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C481").Value = "Metropolitana"
$ws.Range("D481").Value = 44516
$ws.Range("E481").Value = 13
$ws.Range("F481").Value = 100112023
$ws.Range("G481").Value = "Brócoli"
$ws.Range("H481").Value = "Sin especificar"
$ws.Range("I481").Value = "Primera"
$ws.Range("J481").Value = 12700
$ws.Range("K481").Value = 500
$ws.Range("L481").Value = 600
$ws.Range("M481").Value = 554
$ws.Range("N481").Value = "$/unidad"
$ws.Range("O481").Value = "Región Metropolitana"
$ws.Range("P481").Value = 554
$ws.Range("Q481").Value = 1
$ws.Range("R481").Value = "Hortaliza"

# Row 482: Brócoli, Segunda
$ws.Range("A482").Value = 6
$ws.Range("B482").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C482").Value = "Metropolitana"
$ws.Range("D482").Value = 44516
$ws.Range("E482").Value = 13
$ws.Range("F482").Value = 100112023
$ws.Range("G482").Value = "Brócoli"
$ws.Range("H482").Value = "Sin especificar"
$ws.Range("I482").Value = "Segunda"
$ws.Range("J482").Value = 5200
$ws.Range("K482").Value = 400
$ws.Range("L482").Value = 500
$ws.Range("M482").Value = 454
$ws.Range("N482").Value = "$/unidad"
$ws.Range("O482").Value = "Región Metropolitana"
$ws.Range("P482").Value = 454
$ws.Range("Q482").Value = 1
$ws.Range("R482").Value = "Hortaliza"
